$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2667
$ws1.Range("F11").Value = 9854
$ws1.Range("F15").Value = 610
$ws1.Range("F16").Value = 11726
$ws1.Range("F17").Value = 12057

# Sheet "全部类型" (sheet4): same updates, rows shifted by one
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2667
$ws4.Range("F12").Value = 9854
$ws4.Range("F16").Value = 610
$ws4.Range("F17").Value = 11726
$ws4.Range("F18").Value = 12057
